# Delete indifferent experiments and results.
# Updates the registration-quality results for both the liver and tumor
# sheets (re-run with a different rigid-registration metric), renames the
# "01 Rigid MI" / "02 B-spline MI" experiment columns, restyles the
# "Mean" row on the tumor sheet to match the liver sheet, resizes a
# couple of columns, and switches the active tab to the tumor sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("liver.nii.gz")
$ws2 = $wb.Worksheets.Item("tumor.nii.gz")

# ---------------------------------------------------------------------
# Rename the experiment headers (shared by both sheets) in place.
# ---------------------------------------------------------------------
foreach ($ws in @($ws1, $ws2)) {
    $ws.Range("C1").Value = "01 Rigid KS"
    $ws.Range("D1").Value = "02 B-Spline MI"
}

# ---------------------------------------------------------------------
# liver.nii.gz (sheet1) updated Dice scores for columns C (01 Rigid KS)
# and D (02 B-Spline MI).
# ---------------------------------------------------------------------
$ws1.Range("C2").Value  = 0.78769
$ws1.Range("D2").Value  = 0.8675
$ws1.Range("C3").Value  = 0.83594
$ws1.Range("D3").Value  = 0.8858
$ws1.Range("C4").Value  = 0.79145
$ws1.Range("D4").Value  = 0.88726
$ws1.Range("C5").Value  = 0.90461
$ws1.Range("D5").Value  = 0.93025
$ws1.Range("C6").Value  = 0.78029
$ws1.Range("D6").Value  = 0.88214
$ws1.Range("C7").Value  = 0.74066
$ws1.Range("D7").Value  = 0.85975
$ws1.Range("C8").Value  = 0.72307
$ws1.Range("D8").Value  = 0.87766
$ws1.Range("C9").Value  = 0.67135
$ws1.Range("D9").Value  = 0.84387
$ws1.Range("C10").Value = 0.7999
$ws1.Range("D10").Value = 0.84028
$ws1.Range("C11").Value = 0.79255
$ws1.Range("D11").Value = 0.87685
$ws1.Range("C12").Value = 0.67135
$ws1.Range("D12").Value = 0.84028
$ws1.Range("C13").Value = 0.90461
$ws1.Range("D13").Value = 0.93025
$ws1.Range("C14").Value = 0.7836225
$ws1.Range("D14").Value = 0.8768241666666666
$ws1.Range("C15").Value = 0.78769
$ws1.Range("D15").Value = 0.87685

# ---------------------------------------------------------------------
# tumor.nii.gz (sheet2) updated Dice scores for columns C and D.
# ---------------------------------------------------------------------
$ws2.Range("C2").Value  = 0.49798
$ws2.Range("D2").Value  = 0.74748
$ws2.Range("C3").Value  = 0.86615
$ws2.Range("D3").Value  = 0.90085
$ws2.Range("C4").Value  = 0.70783
$ws2.Range("D4").Value  = 0.71076
$ws2.Range("C5").Value  = 0.88663
$ws2.Range("D5").Value  = 0.90903
$ws2.Range("C6").Value  = 0.76722
$ws2.Range("D6").Value  = 0.79607
$ws2.Range("C7").Value  = 0
$ws2.Range("D7").Value  = 0.05171
$ws2.Range("C8").Value  = 0.78838
$ws2.Range("D8").Value  = 0.77805
$ws2.Range("C9").Value  = 0.74794
$ws2.Range("D9").Value  = 0.7337
$ws2.Range("C10").Value = 0.71801
$ws2.Range("D10").Value = 0.7706
$ws2.Range("C11").Value = 0.6164
$ws2.Range("D11").Value = 0.63222
$ws2.Range("D12").Value = 0.05171
$ws2.Range("C13").Value = 0.88663
$ws2.Range("D13").Value = 0.90903
$ws2.Range("C14").Value = 0.6235975
$ws2.Range("D14").Value = 0.6659341666666667
$ws2.Range("C15").Value = 0.71801
$ws2.Range("D15").Value = 0.74748

# ---------------------------------------------------------------------
# Mirror the "Mean" row highlight style (bold font + yellow fill) from
# the liver sheet onto the tumor sheet's Mean row (C14:D14).
# ---------------------------------------------------------------------
$meanRange = $ws2.Range("C14:D14")
$meanRange.Font.Bold = $true
$meanRange.Interior.Pattern = 1
$meanRange.Interior.Color = 65535

# ---------------------------------------------------------------------
# Column width tweaks on the liver sheet (C/D slightly resized) and add
# matching column widths to the tumor sheet.
# ---------------------------------------------------------------------
$ws1.Columns.Item(3).ColumnWidth = 9.666666666666666
$ws1.Columns.Item(4).ColumnWidth = 12.833333333333334

$ws2.Columns.Item(1).ColumnWidth = 19.666666666666668
$ws2.Columns.Item(2).ColumnWidth = 9.5
$ws2.Columns.Item(3).ColumnWidth = 9.666666666666666
$ws2.Columns.Item(4).ColumnWidth = 12.833333333333334

# ---------------------------------------------------------------------
# Selection / active tab: the tumor sheet becomes the active tab. Both
# sheets end up with a (Ctrl-click) multi-selection spanning C14 and
# D14 -- land on the cell that is the true "active" cell of each
# multi-selection (C14 on the liver sheet, D14 on the tumor sheet).
# ---------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("D14").Select()
$ws1.Range("C14").Select()

$ws2.Activate()
$ws2.Range("C14").Select()
$ws2.Range("D14").Select()
